$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the sheet view so column E becomes the top-left visible cell (topLeftCell="E1")
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1

# Adjust column widths (stored as Excel "character" width units in the XML)
$ws.Columns.Item(3).ColumnWidth = 16.140625
$ws.Columns.Item(4).ColumnWidth = 18.28515625
$ws.Columns.Item(10).ColumnWidth = 24.7109375
$ws.Columns.Item(11).ColumnWidth = 34.140625
$ws.Columns.Item(12).ColumnWidth = 43
